$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item('展览')
$ws.Range("F2").Value = 834
$ws.Range("C3").Value = '广州·代号鸢only2.0'
$ws.Range("D3").Value = '清河东路288号 科尔海悦酒店'
$ws.Range("E3").Value = '2024.03.16 10:00-03.16 21:00'
$ws.Range("F3").Value = 1412
$ws.Range("G3").Value = 39
$ws.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=79828'
$ws.Range("I3").Value = '//i0.hdslb.com/bfs/openplatform/202312/RVUVc8oy1702549585918.jpeg'
$ws.Range("C4").Value = '广州·原神X星穹铁道X绝区零ONLY'
$ws.Range("D4").Value = '洛浦街夏滘西环路1号(厦滘地铁站A口步行290米) 厦喾岭南电商园会展中心'
$ws.Range("E4").Value = '2024.03.16 10:00-03.16 17:00'
$ws.Range("F4").Value = 852
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=80715'
$ws.Range("I4").Value = '//i0.hdslb.com/bfs/openplatform/202401/Lt6ZYvA41704878219924.jpeg'
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = '2024-03-23'
$ws.Range("C5").Value = '广州·BanG Dream ONLY'
$ws.Range("D5").Value = '西环路1号 广州岭南会展中心'
$ws.Range("E5").Value = '2024.03.23 10:00-03.23 17:00'
$ws.Range("F5").Value = 485
$ws.Range("G5").Value = 65
$ws.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=81754'
$ws.Range("I5").Value = '//i1.hdslb.com/bfs/openplatform/202402/CtAZIgth1709176182850.jpeg'
$ws.Range("C6").Value = '广州·YU 7th动漫嘉年华'
$ws.Range("D6").Value = '珠江西路8号 高德置地夏广场'
$ws.Range("E6").Value = '2024.03.23 10:00-03.24 17:00'
$ws.Range("F6").Value = 198
$ws.Range("G6").Value = 55
$ws.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=81627'
$ws.Range("I6").Value = '//i1.hdslb.com/bfs/openplatform/202403/kiGKagY41710141580683.jpeg'
$ws.Range("F7").Value = 635
$ws.Range("F8").Value = 198
$ws.Range("F10").Value = 53
$ws.Range("F12").Value = 124
$ws.Range("F13").Value = 1655
$ws.Range("F14").Value = 212
$ws.Range("F15").Value = 36
$ws.Range("F17").Value = 77
$ws.Range("F18").Value = 398
$ws.Range("F20").Value = 637
$ws.Range("F21").Value = 33
$ws.Range("F23").Value = 935
$ws.Range("F25").Value = 1481
$ws.Range("F26").Value = 199

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item('演出')
$ws.Range("F5").Value = 203
$ws.Range("F7").Value = 275
$ws.Range("F8").Value = 65
$ws.Range("F11").Value = 128

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range("F3").Value = 834
$ws.Range("C4").Value = '广州·代号鸢only2.0'
$ws.Range("D4").Value = '清河东路288号 科尔海悦酒店'
$ws.Range("E4").Value = '2024.03.16 10:00-03.16 21:00'
$ws.Range("F4").Value = 1412
$ws.Range("G4").Value = 39
$ws.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=79828'
$ws.Range("I4").Value = '//i0.hdslb.com/bfs/openplatform/202312/RVUVc8oy1702549585918.jpeg'
$ws.Range("C5").Value = '广州·原神X星穹铁道X绝区零ONLY'
$ws.Range("D5").Value = '洛浦街夏滘西环路1号(厦滘地铁站A口步行290米) 厦喾岭南电商园会展中心'
$ws.Range("E5").Value = '2024.03.16 10:00-03.16 17:00'
$ws.Range("F5").Value = 852
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=80715'
$ws.Range("I5").Value = '//i0.hdslb.com/bfs/openplatform/202401/Lt6ZYvA41704878219924.jpeg'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '2024-03-17'
$ws.Range("C6").Value = '广州·三月的幻想演唱会2024「飞越蓝色时刻」'
$ws.Range("D6").Value = '恩宁路265号三层、四层自编01 MAO Livehouse广州(永庆坊店)'
$ws.Range("E6").Value = '2024.03.17 19:00-03.17 20:30'
$ws.Range("F6").Value = 117
$ws.Range("G6").Value = 380
$ws.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=80870'
$ws.Range("I6").Value = '//i1.hdslb.com/bfs/openplatform/202401/8WBT7H6W1705376580145.png'
$ws.Range("C7").Value = '广州·梁祝 ·黄河经典名曲大型管弦交响音乐会'
$ws.Range("D7").Value = '东风中路299号 广州中山纪念堂'
$ws.Range("E7").Value = '2024.03.17 19:30-03.17 21:00'
$ws.Range("F7").Value = 31
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=81788'
$ws.Range("I7").Value = '//i2.hdslb.com/bfs/openplatform/202402/54YX2MVU1707208994883.jpeg'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '2024-03-23'
$ws.Range("C8").Value = '广州·BanG Dream ONLY'
$ws.Range("D8").Value = '西环路1号 广州岭南会展中心'
$ws.Range("E8").Value = '2024.03.23 10:00-03.23 17:00'
$ws.Range("F8").Value = 485
$ws.Range("G8").Value = 65
$ws.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=81754'
$ws.Range("I8").Value = '//i1.hdslb.com/bfs/openplatform/202402/CtAZIgth1709176182850.jpeg'
$ws.Range("C9").Value = '广州·YU 7th动漫嘉年华'
$ws.Range("D9").Value = '珠江西路8号 高德置地夏广场'
$ws.Range("E9").Value = '2024.03.23 10:00-03.24 17:00'
$ws.Range("F9").Value = 198
$ws.Range("G9").Value = 55
$ws.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=81627'
$ws.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202403/kiGKagY41710141580683.jpeg'
$ws.Range("F10").Value = 635
$ws.Range("F12").Value = 198
$ws.Range("F14").Value = 53
$ws.Range("F16").Value = 124
$ws.Range("F17").Value = 1655
$ws.Range("F18").Value = 203
$ws.Range("F19").Value = 212
$ws.Range("F20").Value = 36
$ws.Range("F22").Value = 77
$ws.Range("F23").Value = 398
$ws.Range("F26").Value = 275
$ws.Range("F27").Value = 65
$ws.Range("F28").Value = 637
$ws.Range("F31").Value = 128
$ws.Range("F32").Value = 128
$ws.Range("F33").Value = 33
$ws.Range("F35").Value = 935
$ws.Range("F37").Value = 1481
$ws.Range("F38").Value = 199
